# Reflection.docx edit script
# 1. "Metadatlarin" -> "Metadatalarin" typo fix
# 2. "Assembly partial class" -> "Assembly partial classi" wording fix
# 3. Insert a new "*BaseType()" section before the "*GetConstructor & GetConstructors" section
# 4. Remove 5 superfluous blank paragraphs that used to sit before "*GetField & GetFields"
#    (their "budget" moved up into the new BaseType() section)
# 5. "BindingFields" -> "BindingFlags" typo fix, only in the GetField(...) sample line

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $exactText) {
            return $i
        }
    }
    return -1
}

# --- 1. Metadatlarin -> Metadatalarin --------------------------------------
$null = $d.Content.Find.Execute(
    "Metadatların",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Metadataların", 2)

# --- 2. partial class -> partial classi -------------------------------------
$null = $d.Content.Find.Execute(
    "Assembly partial class vasitəsi",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Assembly partial classı vasitəsi", 2)

# --- 3. Insert the new *BaseType() section -----------------------------------
$anchorIdx = Find-ParagraphIndex $d "*GetConstructor & GetConstructors"
$prev = $d.Paragraphs($anchorIdx - 1)   # blank, non-bold paragraph right above it

# Create the 5 new (still blank / non-bold) paragraphs first, chaining off a
# non-bold paragraph so they all inherit "no bold" cleanly.
$prev.Range.InsertParagraphAfter()
$pHeader = $prev.Next()

$pHeader.Range.InsertParagraphAfter()
$pExplain = $pHeader.Next()

$pExplain.Range.InsertParagraphAfter()
$pCode = $pExplain.Next()

$pCode.Range.InsertParagraphAfter()
$pBlank1 = $pCode.Next()

$pBlank1.Range.InsertParagraphAfter()
$pBlank2 = $pBlank1.Next()

# Fill in the plain-text (non-bold) paragraphs.
$pExplain.Range.InsertAfter("Əldə olunan Type object vasitəsi ilə onun miras aldığı class’ı(type’ı) əldə etmək olar.")
$pCode.Range.InsertAfter("Type? baseType = type.BaseType;")

# Bold only the section header.
$pHeader.Range.Font.Bold = 1
$pHeader.Range.InsertAfter("*BaseType()")

# --- 4. Drop 5 blank paragraphs before "*GetField & GetFields" ---------------
$fieldsIdx = Find-ParagraphIndex $d "*GetField & GetFields"
$delStart = $d.Paragraphs($fieldsIdx - 8).Range.Start
$delEnd = $d.Paragraphs($fieldsIdx - 4).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

# --- 5. BindingFields -> BindingFlags (GetField sample only) -----------------
$null = $d.Content.Find.Execute(
    "type.GetField(fieldName, BindingFields.Public | BindingFields.NonPublic | BindingFields.Instance);",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "type.GetField(fieldName, BindingFlags.Public | BindingFlags.NonPublic | BindingFlags.Instance);", 2)

Write-Output "done"
